# Re-ran the averaged-intensities code, now including three new spiral
# sampling schemes ("Spiral-90deg-10rot-5space", "Spiral-90deg-15rot-5space",
# "Spiral-90deg-10rot-3space"). The Gaussian-Quadrature scheme row moves up
# to sit right after "Ring Perpendicular to TD", the three new spiral rows
# follow it, and the remaining schemes shift down accordingly, growing the
# table from A1:M16 to A1:M19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new index-column cells (A10:A19) the same formatting (bold,
# centered/top, thin border) as the rest of the index column before filling
# in their values.
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Full target content for rows 10-19 (HKL index, scheme name, and the 11
# averaged-intensity columns C:M) after re-running the analysis with the
# spiral schemes included.
$rowsData = @(
    @(8,  "Gaussian-Quadrature",          @(1.029151185967918, 0.8894746492138162, 1.015970039646528, 1.029151185967918, 0.9386335319368188, 1.050516769224966, 1.020595932049202, 0.8894746492138162, 0.9527223444301722, 0.990936765199045, 0.9907236846732083)),
    @(9,  "Spiral-90deg-10rot-5space",    @(0.9865264967460508, 0.947273918617817, 1.014922471336538, 0.9865264967460508, 0.9583465575639657, 1.05653774535591, 1.006027579896071, 0.947273918617817, 0.9810981949771774, 0.9838123458616141, 0.9949391282527253)),
    @(10, "Spiral-90deg-15rot-5space",    @(0.9863314240320836, 0.9482942795311028, 1.014631289389101, 0.9863314240320836, 0.958948227160391, 1.055721081218991, 1.005737730051075, 0.9482942795311028, 0.9814627844601018, 0.9838971042460927, 0.9949440052304572)),
    @(11, "Spiral-90deg-10rot-3space",    @(0.9864858207631485, 0.947476916137293, 1.014876820031628, 0.9864858207631485, 0.9584658806729444, 1.056384977097646, 1.005954059929437, 0.947476916137293, 0.9811768680844608, 0.9838313444238046, 0.9949407457720162)),
    @(12, "NoRotation-tilt60deg",         @(1.012944, 0.8111200000000001, 1.052952, 1.012944, 0.8784680000000011, 1.167387999999997, 1.043996000000001, 0.8111200000000001, 0.9320360000000002, 0.9724900000000003, 0.9944779999999999)),
    @(13, "Rotation-NoTilt",              @(1.05, 0.61, 1.11, 1.05, 0.76, 1.33, 1.1, 0.61, 0.8600000000000001, 0.9550000000000001, 0.9933333333333335)),
    @(14, "Rotation-60detTilt",           @(1.027709146521598, 0.769980440166402, 1.061427330867196, 1.027709146521598, 0.858271062220804, 1.188010268057598, 1.054819286835198, 0.769980440166402, 0.9157038855167989, 0.9717065160191984, 0.9933695891114661)),
    @(15, "HexGrid-90degTilt5degRes",     @(0.993435898491369, 0.9947004371482534, 0.9943060068628734, 0.993435898491369, 0.9932646349701468, 0.9952449745221986, 0.9942550969859905, 0.9947004371482534, 0.9945032220055634, 0.9939695602484663, 0.9942011748301387)),
    @(16, "HexGrid-90degTilt22p5degRes",  @(0.9922415012822449, 1.007309691688417, 0.990682121422112, 0.9922415012822449, 0.9997303942744054, 0.9831589670918084, 0.9900215253719409, 1.007309691688417, 0.9989959065552645, 0.9956187039187547, 0.9938573668551549)),
    @(17, "HexGrid-60degTilt5degRes",     @(0.9851567748163171, 1.038822212636244, 0.9823845624334122, 0.9851567748163171, 1.01941278468194, 0.9613684355296636, 0.9816636803529625, 1.038822212636244, 1.010603387534828, 0.9978800811755727, 0.9948014084084232))
)

$row = 10
foreach ($entry in $rowsData) {
    $hklIndex = $entry[0]
    $schemeName = $entry[1]
    $values = $entry[2]

    $ws.Cells.Item($row, 1).Value2 = $hklIndex
    $ws.Cells.Item($row, 2).Value = $schemeName

    for ($col = 0; $col -lt $values.Count; $col++) {
        $ws.Cells.Item($row, 3 + $col).Value2 = $values[$col]
    }

    $row = $row + 1
}

"updated rows 10-19 with spiral-scheme data"
